# Saldo_guide.xlsx update: refresh the "Dt. Referencia" (column G) extract
# date from 2024-10-17 to 2024-10-18 for every data row, rename the sheet
# to match the new extract timestamp, and correct a handful of balance
# figures (columns D/E/H) that were re-extracted for this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename sheet to reflect the new export timestamp ---------------------
$ws.Name = "IClientBalance-20241018-090426-"

# --- Bump every reference date in column G (rows 2-274) from 45582 --------
# --- (2024-10-17) to 45583 (2024-10-18). A scalar assignment to the whole
# --- range broadcasts the value to every cell while preserving formatting.
$ws.Range("G2:G274").Value = 45583

# --- Row-specific balance corrections --------------------------------------
# Row 15: Saldo Previsto / Vl. Total re-extracted
$ws.Range("E15").Value = 286.22000000000003
$ws.Range("H15").Value = 286.22000000000003

# Row 51: Saldo Previsto / Vl. Total re-extracted
$ws.Range("E51").Value = 46081.98
$ws.Range("H51").Value = 46081.98

# Row 104: Saldo Previsto / Vl. Total re-extracted
$ws.Range("E104").Value = 269.5
$ws.Range("H104").Value = 269.5

# Row 108: Vl. Projetado now populated; Vl. Total = Vl. Projetado + Saldo Previsto
$ws.Range("D108").Value = 90428.04
$ws.Range("H108").Value = 91271.95

# Row 118: Saldo Previsto / Vl. Total re-extracted
$ws.Range("E118").Value = 16828.830000000002
$ws.Range("H118").Value = 16828.830000000002

# Row 189: Saldo Previsto / Vl. Total re-extracted
$ws.Range("E189").Value = 43.79
$ws.Range("H189").Value = 43.79

# Row 232: Saldo Previsto / Vl. Total re-extracted
$ws.Range("E232").Value = 12016.61
$ws.Range("H232").Value = 12016.61
